# Updated cryptos list values (price + 1h volume change) per upstream diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.678.28"
$ws.Range("E2").Value = "  +3.11%  "
$ws.Range("D3").Value = "3.002.10"
$ws.Range("E3").Value = "  +2.85%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'510.24"
$ws.Range("E5").Value = "  +7.50%  "
$ws.Range("D6").Value = "'139.70"
$ws.Range("E6").Value = "  +8.72%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'0.433"
$ws.Range("E8").Value = "  +5.71%  "
$ws.Range("E9").Value = "  +12.58%  "
$ws.Range("E10").Value = "  +10.55%  "
$ws.Range("E11").Value = "  +4.63%  "
$ws.Range("E12").Value = "  +4.26%  "
$ws.Range("D13").Value = "3.516.89"
$ws.Range("E13").Value = "  +2.84%  "
$ws.Range("D14").Value = "'25.70"
$ws.Range("E14").Value = "  +8.72%  "
$ws.Range("D15").Value = "'0.0000155"
$ws.Range("E15").Value = "  +14.42%  "
$ws.Range("D16").Value = "56.753.67"
$ws.Range("E16").Value = "  +3.46%  "
$ws.Range("D17").Value = "3.000.11"
$ws.Range("E17").Value = "  +2.97%  "
$ws.Range("E18").Value = "  +8.64%  "
$ws.Range("D19").Value = "'12.48"
$ws.Range("E19").Value = "  +7.38%  "
$ws.Range("E20").Value = "  +9.22%  "
$ws.Range("D21").Value = "'329.39"
$ws.Range("E21").Value = "  +7.97%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  -0.24%  "
$ws.Range("E23").Value = "  +7.91%  "
$ws.Range("D24").Value = "'62.75"
$ws.Range("E24").Value = "  +5.31%  "
$ws.Range("D25").Value = "'0.173"
$ws.Range("E25").Value = "  +12.75%  "
$ws.Range("E26").Value = "  +0.71%  "
$ws.Range("D27").Value = "0.0₃0912"
$ws.Range("E27").Value = "  +11.28%  "
$ws.Range("D28").Value = "'6.71"
$ws.Range("E28").Value = "  +7.28%  "
$ws.Range("D29").Value = "'7.10"
$ws.Range("E29").Value = "  +12.17%  "
$ws.Range("E30").Value = "  +11.30%  "
$ws.Range("E31").Value = "  +8.74%  "
$ws.Range("D32").Value = "'20.67"
$ws.Range("E32").Value = "  +8.93%  "
$ws.Range("D33").Value = "'155.05"
$ws.Range("E33").Value = "  +6.97%  "
$ws.Range("E34").Value = "  +7.46%  "
$ws.Range("E35").Value = "  +3.70%  "
$ws.Range("E36").Value = "  +3.57%  "
$ws.Range("E37").Value = "  +8.88%  "
$ws.Range("D38").Value = "'24.26"
$ws.Range("E38").Value = "  +2.74%  "
$ws.Range("D39").Value = "3.035.61"
$ws.Range("E39").Value = "  +3.03%  "
$ws.Range("E40").Value = "  +3.53%  "
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("E42").Value = "  +5.16%  "
$ws.Range("D43").Value = "2.268.83"
$ws.Range("E43").Value = "  +9.90%  "
$ws.Range("E44").Value = "  +5.51%  "
$ws.Range("E45").Value = "  +3.20%  "
$ws.Range("D46").Value = "'3.66"
$ws.Range("E46").Value = "  +5.91%  "
$ws.Range("D47").Value = "'1.99"
$ws.Range("E47").Value = "  +22.55%  "
$ws.Range("E48").Value = "  +9.08%  "
$ws.Range("E49").Value = "  +7.16%  "
$ws.Range("E50").Value = "  +6.84%  "
$ws.Range("D51").Value = "'0.0869"
$ws.Range("E51").Value = "  +8.72%  "
